$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee Census")
$ws.Range("A4:N8").ClearContents()
$ws.Range("G4").Hyperlinks.Delete()
